$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "Colchisin (Colchi)" size changes from "20 cc Tablet" to "500 mg Tablet"
$ws.Range("B3").Value = "500 mg Tablet"

# New row 4: second Colchisin (Colchi) variant - "20 cc Injection"
$ws.Range("A4").Value = "Colchisin (Colchi)"
$ws.Range("B4").Value = "20 cc Injection"
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 1

# Row heights (auto re-fit as observed after edits)
$ws.Rows.Item(3).RowHeight = 14.9
$ws.Rows.Item(4).RowHeight = 13.8

# Selection / view state left at B8 after editing
[void]$ws.Range("B8").Select()
